$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 20 data, mirroring the style/format of the existing rows
$ws.Range("A19").Copy()
$ws.Range("A20").PasteSpecial(-4122)
$ws.Range("A20").Value = 45986

$ws.Range("B20").Value = 2025
$ws.Range("C20").Value = -2.06674933094535
$ws.Range("D20").Value = 2026
$ws.Range("E20").Value = -1.12081074591468
